# ------------------------------------------------------------------
# Adds a new "2022-Q3" sheet (fund-holdings detail) right after the
# "总计" (totals) sheet, and inserts the corresponding summary row at
# the top of the "总计" table, shifting the existing quarters down.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------- 1. Create the new "2022-Q3" detail sheet ---------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet    = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Reuse the existing header / index-column formatting by copying a
# same-shaped block from a sheet that already has >= 22 rows, so the
# new sheet picks up the same shared cell style (bold header row +
# bordered index column) instead of inventing new style records.
$fmtSource = $wb.Worksheets.Item("2021-Q2")
$fmtSource.Range("A1:H22").Copy()
$newSheet.Range("A1:H22").PasteSpecial(-4122)

# Header row (row 1)
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Fund holdings data for 2022-Q3 (21 rows)
$q3Data = @(
    @('166005','中欧价值发现混合 -A','26.62','93.73','3.83','1.0195',9),
    @('001810','中欧潜力价值灵活配置混合A','19.07','93.66','3.44','0.6560',9),
    @('010744','工银灵动价值混合A','11.46','75.15','4.93','0.5650',2),
    @('004232','中欧价值发现混合 -C','8.18','93.73','3.83','0.3133',9),
    @('000574','宝盈新价值灵活配置混合A','6.15','89.60','4.36','0.2681',10),
    @('003715','宝盈消费主题灵活配置混合','4.42','89.85','4.47','0.1976',9),
    @('166024','中欧恒利三年定期开放混合','3.99','98.45','3.94','0.1572',6),
    @('009223','宝盈现代服务业混合A','3.18','93.59','4.58','0.1456',9),
    @('001651','工银新蓝筹股票A','3.65','82.09','3.31','0.1208',7),
    @('481013','工银消费服务混合A','2.76','70.46','3.52','0.0972',2),
    @('005764','中欧潜力价值灵活配置混合C','2.01','93.66','3.44','0.0691',9),
    @('010745','工银灵动价值混合C','0.87','75.15','4.93','0.0429',2),
    @('008303','宝盈龙头优选股票A','0.65','93.98','4.57','0.0297',10),
    @('007574','宝盈新价值灵活配置混合C','0.66','89.60','4.36','0.0288',10),
    @('009224','宝盈现代服务业混合C','0.43','93.59','4.58','0.0197',9),
    @('001882','中欧价值发现混合 -E','0.43','93.73','3.83','0.0165',9),
    @('006675','宝盈品牌消费股票A','0.17','93.83','4.56','0.0078',10),
    @('008304','宝盈龙头优选股票C','0.16','93.98','4.57','0.0073',10),
    @('006676','宝盈品牌消费股票C','0.14','93.83','4.56','0.0064',10),
    @('011476','工银新蓝筹股票C','0.12','82.09','3.31','0.0040',7),
    @('011475','工银消费服务混合C','0.02','70.46','3.52','0.0007',2)
)

# Columns B..G must stay TEXT (fund codes keep leading zeros, the
# percentages/amounts keep their original decimal formatting) - force
# the whole block to Text format before writing, same as source data.
$newSheet.Range("B2:G22").NumberFormat = "@"

$r = 2
foreach ($row in $q3Data) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------- 2. Update the "总计" summary table ----------------------
# Grow the formatted index-column/border down into the new last row
# (row 9) by duplicating row 8 (its style) one row down, then
# overwrite every data row with the final values (including the new
# 2022-Q3 entry at the top).
$totalSheet.Rows.Item(8).Copy()
$totalSheet.Rows.Item(9).Insert()

$totalData = @(
    @(0, "2022-Q3", 21, 3.77),
    @(1, "2022-Q2", 12, 4.16),
    @(2, "2022-Q1", 15, 6.43),
    @(3, "2021-Q4", 18, 7.42),
    @(4, "2021-Q3", 22, 5.01),
    @(5, "2021-Q2", 23, 4.27),
    @(6, "2021-Q1", 11, 4.41),
    @(7, "2020-Q4", 9, 4.4)
)

$r = 2
foreach ($row in $totalData) {
    $totalSheet.Cells.Item($r, 1).Value = $row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = $row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
    $r++
}
